$wb = $excel.ActiveWorkbook
$dash = $wb.Worksheets.Item("dashboard")
$ws = $wb.Worksheets.Item("writing")

# --- Append a new day of progress data (row 59) ---
# Copy row 58's formatting (incl. the date number format on column A) down to row 59.
$ws.Range("A58").Copy($ws.Range("A59"))

$ws.Range("A59").Value = 44225
$ws.Range("B59").Value = 631
$ws.Range("C59").Value = 260
$ws.Range("D59").Value = 87
$ws.Range("E59").Value = 229
$ws.Range("F59").Value = 492
$ws.Range("G59").Value = 416
$ws.Range("H59").Value = 8656
$ws.Range("I59").Value = 19908
$ws.Range("J59").Value = 6457
$ws.Range("K59").Value = 8497
$ws.Range("L59").Value = 186
$ws.Range("M59").Value = 496
$ws.Range("N59").Value = 3575
$ws.Range("O59").Value = 1329
$ws.Range("P59").Formula = "=SUM(C59:O59)"
# Array-enter the "Daily" delta formula (matches the CSE array formula used by the rest of column Q)
# while the cell is still outside the table's boundary (ListObjects reject FormulaArray edits).
$ws.Range("Q59").FormulaArray = "=SUM(ABS(C59:O59-C58:O58))"

# Register the new row with the worksheet's table so ref/autoFilter grow from Q58 to Q59.
$lo = $ws.ListObjects.Item(1)
$newListRow = $lo.ListRows.Add()

# --- Extend the line-chart series to include the new date ---
$chart = $dash.ChartObjects(1).Chart
$dailySeries = $chart.SeriesCollection(1)
$dailySeries.XValues = "=writing!`$A`$2:`$A`$59"
$dailySeries.Values = "=writing!`$Q`$2:`$Q`$59"
$totalSeries = $chart.SeriesCollection(2)
$totalSeries.XValues = "=writing!`$A`$2:`$A`$59"
$totalSeries.Values = "=writing!`$P`$2:`$P`$59"

# --- Switch the active tab from "dashboard" to "writing" ---
[void]$ws.Activate()
[void]$ws.Range("N56").Select()
